$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) for rows 2-6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1049
$ws1.Range("F3").Value = 311
$ws1.Range("F4").Value = 2830
$ws1.Range("F5").Value = 66
$ws1.Range("F6").Value = 598

# Sheet "全部类型" - update column F (想去人数) for rows 4-8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1049
$ws4.Range("F5").Value = 311
$ws4.Range("F6").Value = 2830
$ws4.Range("F7").Value = 66
$ws4.Range("F8").Value = 598
